$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shift the existing "tail" block (old rows 16-62) down by 9 rows so it
#    becomes rows 25-71, making room to insert 9 new data rows at 16-24.
#    We deliberately avoid Range.Insert()/Rows.Insert() here: in this
#    environment inserting a row right below a formatted block causes Excel
#    to fabricate extra (unused) cell style records in styles.xml. Doing the
#    shift manually with Copy/PasteSpecial keeps styles.xml untouched.
# ---------------------------------------------------------------------------
for ($r = 62; $r -ge 16; $r--) {
    $destRow = $r + 9
    $ws.Rows.Item($r).Copy()
    $ws.Rows.Item($destRow).PasteSpecial(-4104) | Out-Null   # xlPasteAll
}
$excel.CutCopyMode = $false

# Make sure every shifted row keeps/gets its correct row height, even the
# ones that hold no cell values at all (Copy/PasteSpecial of a valueless
# row does not always force the row record to materialize on save).
$tailHeights = @{
    25=18.75; 26=18.75; 27=18.75; 28=18.75; 29=18.75; 30=18.75; 31=18.75;
    32=18.75; 33=18.75; 34=18.75; 35=18.75; 36=18.75; 37=18.75; 38=18.75;
    39=18.75; 40=18.75; 41=18.75; 42=18.75; 43=18.75; 44=18.75; 45=18.75;
    46=18.75; 47=18.75; 48=18.75; 49=18.75; 50=18.75; 51=18.75; 52=18.75;
    53=18.75; 54=18.75; 55=31.5;  56=18.75; 57=18.75; 58=18.75; 59=18.75;
    60=18.75; 61=18.75; 62=18.75; 63=18.75; 64=15.65; 65=18.75; 66=18.75;
    67=18.75; 68=18.75; 69=18.75; 70=18.75; 71=28.25
}
foreach ($r in $tailHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $tailHeights[$r]
}

# ---------------------------------------------------------------------------
# 2) Build the 9 new data rows (16-24) using the same look as rows 13-15:
#    style s=4/5/6/5/5 on columns B/C/D/E/F and a 24pt row height (row 18
#    ends up a touch shorter, at 23pt).
# ---------------------------------------------------------------------------
for ($r = 16; $r -le 24; $r++) {
    $ws.Range("B$r").ClearContents()
    $ws.Range("C$r").ClearContents()
    $ws.Range("D$r").ClearContents()
    $ws.Range("E$r").ClearContents()
    $ws.Range("F$r").ClearContents()

    $ws.Range("B15").Copy()
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("C15").Copy()
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("D15").Copy()
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("E15").Copy()
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("F15").Copy()
    $ws.Range("F$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Rows.Item($r).RowHeight = 24
}
$ws.Rows.Item(18).RowHeight = 23

# ---------------------------------------------------------------------------
# 3) Fill in the new course content.
#    Day8 row (12) gets the "static keyword" topic + its recording link;
#    a brand-new Day9 row (13) gets the OOP / class / constructor topic.
#    Values are written in this order so the newly created shared-string
#    table entries land in the same order as the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = "https://youtu.be/K9Rvor70Aiw"
$ws.Range("B13").Value = "day9"
$ws.Range("C13").Value = "oops start and learn about class and objects ,methodoverloading,consturctor and this keyword in java"
$ws.Range("C12").Value = "lean about static keyword in details"

# ---------------------------------------------------------------------------
# 4) Update the active selection to C12, matching the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("C12").Select()
